$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '28.675.05'
$ws.Range("E2").Value = '  -2.11%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.797.15'
$ws.Range("E3").Value = '  -1.86%  '
$ws.Range("E4").Value = '  +0.03%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '231.61'
$ws.Range("E5").Value = '  -1.70%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.5880'
$ws.Range("E6").Value = '  -2.55%  '
$ws.Range("E7").Value = '  -0.02%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.2769'
$ws.Range("E8").Value = '  -0.80%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.06762'
$ws.Range("E9").Value = '  -4.22%  '
$ws.Range("E10").Value = '  -1.46%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.07534'
$ws.Range("E11").Value = '  -1.64%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '1.796.25'
$ws.Range("E12").Value = '  -1.77%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '4.787'
$ws.Range("E13").Value = '  -0.33%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.6160'
$ws.Range("E14").Value = '  -2.12%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '2.041.36'
$ws.Range("E15").Value = '  -1.87%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.000009098'
$ws.Range("E16").Value = '  -8.32%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '75.32'
$ws.Range("E17").Value = '  -4.93%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '28.664.30'
$ws.Range("E18").Value = '  -1.93%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '5.472'
$ws.Range("E19").Value = '  -6.49%  '
$ws.Range("E20").Value = '  -0.07%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '210.23'
$ws.Range("E21").Value = '  -6.44%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '11.49'
$ws.Range("E22").Value = '  -1.96%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '6.826'
$ws.Range("E23").Value = '  -2.89%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '1.005'
$ws.Range("E24").Value = '  +0.01%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '153.52'
$ws.Range("E25").Value = '  -1.51%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '8.084'
$ws.Range("E26").Value = '  +1.16%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '0.1261'
$ws.Range("E27").Value = '  -3.71%  '
$ws.Range("E28").Value = '  -1.00%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '1.427'
$ws.Range("E29").Value = '  -3.60%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '0.06118'
$ws.Range("E30").Value = '  -3.99%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '1.422'
$ws.Range("E31").Value = '  -1.87%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '3.803'
$ws.Range("E32").Value = '  -0.12%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '3.786'
$ws.Range("E33").Value = '  -1.65%  '
$ws.Range("E34").Value = '  -0.08%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.047'
$ws.Range("E35").Value = '  -5.70%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.6416'
$ws.Range("E36").Value = '  -0.84%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '2.501'
$ws.Range("E37").Value = '  -1.87%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '2.716'
$ws.Range("E38").Value = '  -1.04%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '6.442'
$ws.Range("E39").Value = '  -1.29%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.01699'
$ws.Range("E40").Value = '  -2.75%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '1.141.80'
$ws.Range("E41").Value = '  -6.42%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.8828'
$ws.Range("E42").Value = '  -1.89%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '1.006'
$ws.Range("E43").Value = '  +0.21%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '100.10'
$ws.Range("E44").Value = '  -0.37%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '1.949.22'
$ws.Range("E45").Value = '  -2.20%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '60.05'
$ws.Range("E46").Value = '  -4.23%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.00000000111'
$ws.Range("E47").Value = '  -4.55%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '1.587'
$ws.Range("E48").Value = '  +0.47%  '
$ws.Range("B49").Value = 'EnergySwap'
$ws.Range("C49").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '8.384'
$ws.Range("E49").Value = '  -2.28%  '
$ws.Range("B50").Value = 'Cronos'
$ws.Range("C50").Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.05489'
$ws.Range("E50").Value = '  -0.32%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.4482'
$ws.Range("E51").Value = '  -1.64%  '
